# repull data, push all data, mean calculation
# Update dSF (column F) values to match re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 8
    6  = -3
    14 = -2
    15 = 0
    16 = -3
    17 = 4
    20 = -8
    21 = -4
    22 = -1
    24 = -1
    29 = -2
    37 = -5
    38 = 1
    40 = -5
    41 = 0
    42 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
